$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "CS101"
$ws.Range("B2").Value = "Item 2"
$ws.Range("C2").Value = 100
$ws.Range("D2").Value = 999
$ws.Range("E2").Value = "Desc 1"
$ws.Range("F2").Value = "Overlapped"

# Row 3
$ws.Range("A3").Value = "CS103"
$ws.Range("B3").Value = "Item 1"
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = 100
$ws.Range("E3").Value = "Desc 1"
$ws.Range("F3").Value = "Completed"
